# Mark the "butterflies", "country", "herpetology", "macro" and "moody"
# hub rows as added: tick the checkbox in column A and fill in "Y" for
# the feature columns that are now supported (mirrors the formatting of
# the other already-completed rows in the sheet, which use a Text number
# format on the populated cells).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 22 - butterflies
$ws.Range("A22").Value = "[X] "
$ws.Range("C22").Value = "Y"
$ws.Range("C22").NumberFormat = "@"
$ws.Range("D22").Value = "Y"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("G22").Value = "Y"
$ws.Range("G22").NumberFormat = "@"

# Row 34 - country
$ws.Range("A34").Value = "[X] "
$ws.Range("C34").Value = "Y"
$ws.Range("C34").NumberFormat = "@"
$ws.Range("D34").Value = "Y"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("G34").Value = "Y"
$ws.Range("G34").NumberFormat = "@"

# Row 48 - herpetology
$ws.Range("A48").Value = "[X] "
$ws.Range("C48").Value = "Y"
$ws.Range("C48").NumberFormat = "@"
$ws.Range("D48").Value = "Y"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("E48").Value = "Y"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("G48").Value = "Y"
$ws.Range("G48").NumberFormat = "@"

# Row 59 - macro
$ws.Range("A59").Value = "[X] "
$ws.Range("C59").Value = "Y"
$ws.Range("C59").NumberFormat = "@"
$ws.Range("D59").Value = "Y"
$ws.Range("D59").NumberFormat = "@"
$ws.Range("E59").Value = "Y"
$ws.Range("E59").NumberFormat = "@"
$ws.Range("G59").Value = "Y"
$ws.Range("G59").NumberFormat = "@"

# Row 62 - moody
$ws.Range("A62").Value = "[X] "
$ws.Range("C62").Value = "Y"
$ws.Range("C62").NumberFormat = "@"
$ws.Range("D62").Value = "Y"
$ws.Range("D62").NumberFormat = "@"
$ws.Range("E62").Value = "Y"
$ws.Range("E62").NumberFormat = "@"
$ws.Range("G62").Value = "Y"
$ws.Range("G62").NumberFormat = "@"
